$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5 (Ivanov / lecture1_task_2) is removed entirely - the new scraper run
# only produced 3 review rows instead of 4.
# ---------------------------------------------------------------------------
$ws.Rows(5).Delete()

# ---------------------------------------------------------------------------
# Shared "error" comment used by rows 2 and 3 (the new run failed before it
# could read the student's code).
# ---------------------------------------------------------------------------
$errorComment = @'
Error during review: 'RepositoryService' object has no attribute 'read_student_code'
'@

# ---------------------------------------------------------------------------
# Row 2: Petrov -> xesilver, scores reset to 0 because the review errored out.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "xesilver"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = $errorComment
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "2025-09-19 18:02:19"

# ---------------------------------------------------------------------------
# Row 3: Petrov -> xesilver, scores reset to 0 because the review errored out.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "xesilver"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = $errorComment
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "2025-09-19 18:02:19"

# ---------------------------------------------------------------------------
# Row 4: Ivanov -> xesilver, task renamed, new review text/scores from the
# re-run against the freshly downloaded GitHub repo.
# ---------------------------------------------------------------------------
$row4Comment = @'
- What’s working:
  - The calculator handles the four basic operations and prints results.
  - Division by zero is gracefully handled with a message.
  - Basic input prompts are clear and straightforward.
- Areas for improvement (technical correctness and robustness):
  - No input validation beyond division by zero. If a user enters a non-numeric value, float(...) will raise a ValueError and crash the program.
  - Code mixes calculation logic with I/O (top-level input/output). This makes testing harder and reduces reusability.
  - There is no main guard (if __name__ == "__main__":). If this file is imported, it will execute immediately, which is undesirable for a library/module.
  - No unit tests or small testable modules are provided.
- Suggested improvements:
  - Add input validation and exception handling for numeric input.
  - Separate concerns: put calculation logic in pure functions (they already exist) and wrap I/O in a main() function.
  - Add a main guard to prevent running on import.
  - Introduce docstrings for functions to explain behavior and types.
  - Consider using a dispatch map for operations to reduce repetitive conditional logic.
- Example of improved approach (illustrative, not required to copy exactly):
  - Include type hints and docstrings.
  - Use a dict to map operators to functions.
  - Add a simple main() and a guard.
  Example (conceptual):
  def add(x: float, y: float) -> float:
      """Return the sum of x and y."""
      return x + y
  def sub(x: float, y: float) -> float:
      """Return the difference x - y."""
      return x - y
  def mul(x: float, y: float) -> float:
      """Return the product x * y."""
      return x * y
  def div(x: float, y: float) -> float:
      """Return the division x / y."""
      return x / y
  def main():
      try:
          a = float(input("Enter first number: "))
          b = float(input("Enter second number: "))
      except ValueError:
          print("Invalid number")
          return
      op = input("Enter operation (+, -, *, /): ")
      ops = {'+': add, '-': sub, '*': mul, '/': div}
      if op not in ops:
          print("Invalid operation")
          return
      if op == '/' and b == 0:
          print("Cannot divide by zero")
          return
      result = ops[op](a, b)
      print(result)
  if __name__ == "__main__":
      main()
- Why this helps:
  - Improves robustness against bad input.
  - Makes testing easier (you can unit test add, sub, mul, div separately).
  - Improves readability and maintainability with docstrings and explicit flow control.
- Additional notes:
  - Documentation (docstrings) is currently missing. Adding docstrings (and optional type hints) would significantly improve Documentation score.
  - If this is meant as a standalone script, the current approach is acceptable, but adding a __main__ guard and basic input validation would still be beneficial.
'@

$ws.Range("A4").Value = "xesilver"
$ws.Range("C4").Value = "task_1"
$ws.Range("D4").Value = 66
$ws.Range("E4").Value = $row4Comment
$ws.Range("F4").Value = 75
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = "2025-09-19 18:08:39"

# ---------------------------------------------------------------------------
# Header row (A1:J1) loses its bold-white-on-blue look; it becomes plain bold
# black text, no fill, top-aligned instead of centered vertically. Start from
# a clean slate so the blue fill / white font don't linger underneath.
# ---------------------------------------------------------------------------
$header = $ws.Range("A1:J1")
$header.ClearFormats()
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Unfreeze the header row (the frozen pane + its special selection go away).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A1").Select()
